$wb = $excel.ActiveWorkbook

# Map of cell -> new value that must be applied identically to the
# "展览" (index 1) and "全部类型" (index 4) worksheets.
$updates = @{
    "F2"  = 1394
    "F3"  = 2717
    "F4"  = 592
    "F6"  = 6635
    "F7"  = 1020
    "F8"  = 11
    "F9"  = 17
    "F12" = 9
}

$sheetIndexes = @(1, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
